$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Day" column's final row has a distinct date-only number format
# (YYYY-MM-DD) compared to the regular rows above it (YYYY-MM-DD HH:MM:SS).
# Capture both formats before making changes, since the final row's format
# needs to move down to the new last row.
$lastRowDateFormat = $ws.Range("A24").NumberFormat
$normalDateFormat   = $ws.Range("A2").NumberFormat

# Append the new day's data as row 25 (daily update).
$ws.Range("A25").Value = 45974
$ws.Range("B25").Value = 54
$ws.Range("C25").Value = 62
$ws.Range("D25").Value = 61

# Row 25 is now the final row, so it takes on the final-row date format.
$ws.Range("A25").NumberFormat = $lastRowDateFormat

# Row 24 is no longer the final row, so it reverts to the regular date format.
$ws.Range("A24").NumberFormat = $normalDateFormat
